# Rename the sheet and replace the data with a new header-only row:
#   id | discount | name
# (dropping the old order_id/product_id columns and all of the sample
# data rows, as the query output now only needs these three columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab.
$ws.Name = "Query A"

# B1 already carries the bold/bordered/centered header style used for the
# column titles. Copy that formatting onto A1 so the new "id" header cell
# looks the same as the others before we touch any values.
$ws.Range("B1").Copy($ws.Range("A1"))

# Drop all of the old sample data rows (2-5).
$ws.Range("A2:E5").Delete()

# Drop the old D1 ("name") / E1 ("discount") header cells; the remaining
# A1:C1 header cells will be given the final labels below.
$ws.Range("D1:E1").Clear()

# Set the final header labels for the three remaining columns.
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "discount"
$ws.Range("C1").Value = "name"
